$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The edit inserts a new column H ("InDataLower") that holds the lower-cased
# version of what used to be column H ("InData" - mixed-case gene tags).
# The old column H (with its data + styling) shifts right to become column I.
# ---------------------------------------------------------------------------

# Insert a new blank column before the existing column H (8). This pushes
# the old H (data + per-cell style) one column to the right, into I.
$ws.Columns.Item(8).Insert()

# The freshly inserted column inherits formatting from its neighbour; strip
# that so the new column starts out unstyled (matches the target workbook).
$ws.Columns.Item(8).ClearFormats()

# New header for the inserted column.
$ws.Range("H1").Value2 = "InDataLower"

$lastRow = 102
for ($r = 2; $r -le $lastRow; $r++) {
    $oldVal = $ws.Range("I$r").Value2
    if ($oldVal -ne $null -and $oldVal -ne "") {
        $ws.Range("H$r").Value2 = $oldVal.ToLower()
    }
}
